$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.174.03"
$ws.Range("E2").Value = "  +3.21%  "
$ws.Range("D3").Value = "3.814.63"
$ws.Range("E3").Value = "  +1.26%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "708.76"
$ws.Range("E5").Value = "  +12.74%  "
$ws.Range("E6").Value = "  +5.06%  "
$ws.Range("D7").Value = "3.814.00"
$ws.Range("E7").Value = "  +1.25%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").Value = "0.527"
$ws.Range("E9").Value = "  +1.50%  "
$ws.Range("D10").Value = "0.165"
$ws.Range("E10").Value = "  +4.04%  "
$ws.Range("D11").Value = "7.41"
$ws.Range("E11").Value = "  +9.81%  "
$ws.Range("D12").Value = "0.464"
$ws.Range("E12").Value = "  +1.76%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000257"
$ws.Range("E13").Value = "  +7.81%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.50"
$ws.Range("E14").Value = "  +5.04%  "
$ws.Range("D15").Value = "4.454.69"
$ws.Range("E15").Value = "  +1.22%  "
$ws.Range("D16").Value = "3.813.59"
$ws.Range("E16").Value = "  +1.14%  "
$ws.Range("D17").Value = "71.187.82"
$ws.Range("E17").Value = "  +3.26%  "
$ws.Range("D18").Value = "18.02"
$ws.Range("E18").Value = "  +2.21%  "
$ws.Range("E19").Value = "  +3.68%  "
$ws.Range("E20").Value = "  +0.51%  "
$ws.Range("D21").Value = "11.27"
$ws.Range("E21").Value = "  +18.54%  "
$ws.Range("D22").Value = "485.13"
$ws.Range("E22").Value = "  +5.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.720"
$ws.Range("E23").Value = "  +2.53%  "
$ws.Range("D24").Value = "84.12"
$ws.Range("E24").Value = "  +2.45%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000146"
$ws.Range("E25").Value = "  +1.44%  "
$ws.Range("D26").Value = "12.57"
$ws.Range("E26").Value = "  +4.05%  "
$ws.Range("D27").Value = "10.74"
$ws.Range("E27").Value = "  +6.37%  "
$ws.Range("D28").Value = "2.21"
$ws.Range("E28").Value = "  +4.38%  "
$ws.Range("D29").Value = "3.965.80"
$ws.Range("E29").Value = "  +1.23%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  -0.08%  "
$ws.Range("D31").Value = "3.15"
$ws.Range("E31").Value = "  +18.17%  "
$ws.Range("D32").Value = "7.63"
$ws.Range("E32").Value = "  +8.18%  "
$ws.Range("E33").Value = "  +2.09%  "
$ws.Range("D34").Value = "29.74"
$ws.Range("E34").Value = "  +4.87%  "
$ws.Range("D35").Value = "0.179"
$ws.Range("E35").Value = "  +2.59%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.30"
$ws.Range("E36").Value = "  +4.68%  "
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("D38").Value = "3.765.74"
$ws.Range("E38").Value = "  +1.24%  "
$ws.Range("E39").Value = "  +3.75%  "
$ws.Range("E40").Value = "  +7.13%  "
$ws.Range("D41").Value = "6.02"
$ws.Range("E41").Value = "  +4.40%  "
$ws.Range("E42").Value = "  +14.64%  "
$ws.Range("D43").Value = "0.000335"
$ws.Range("E43").Value = "  +27.35%  "
$ws.Range("D44").Value = "0.974"
$ws.Range("E44").Value = "  +1.45%  "
$ws.Range("D45").Value = "0.999"
$ws.Range("E45").Value = "  -0.05%  "
$ws.Range("D47").Value = "45.73"
$ws.Range("E47").Value = "  +6.61%  "
$ws.Range("D48").Value = "49.55"
$ws.Range("E48").Value = "  +5.59%  "
$ws.Range("D49").Value = "160.84"
$ws.Range("E49").Value = "  +2.60%  "
$ws.Range("E50").Value = "  +0.12%  "
$ws.Range("E51").Value = "  +3.27%  "
